# Adding "Area" (col G) and "Atotal" (col H) to the discharge worksheet,
# mirroring the existing velocity-area (Q/Qtotal) computation pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"

# Row 2: area of the first (near-bank) segment measured from 0, and the
# running area total across the first block of segments
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# Row 3: area of the next segment, referencing the previous midpoint
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

# Rows 4-15: the same incremental-area pattern, one formula per row so the
# relative D/B references shift down with the row (same shape as the
# existing shared D/E formulas just below)
for ($r = 4; $r -le 15; $r++) {
    $prev = $r - 1
    $ws.Range("G$r").Formula = "=(D$r-D$prev)*B$r/100"
}

# Leave the selection on H2, matching the sheet's new focal cell
$ws.Range("H2").Select()
